$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number need an
# explicit Text format so they stay strings (matching the source data, which
# stores every Price/Volume cell as text).

$ws.Range("D2").Value = "42.696.61"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.309.59"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.92"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.64"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.35"
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.64"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.981"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "2.661.95"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "2.316.63"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "42.690.82"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +33.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.00"
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.57"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.98"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.78"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.28"
$ws.Range("E30").Value = "  +8.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.03"
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("E32").Value = "  +7.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("E33").Value = "  +2.94%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  -10.11%  "
$ws.Range("E36").Value = "  +0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.63"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0355"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  +9.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.65"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.75"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.95"
$ws.Range("E47").Value = "  +8.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.96"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("E49").Value = "  +1.35%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "1.623.36"
$ws.Range("E51").Value = "  +6.38%  "
